$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 - new Mac-Address entry
$ws.Range("A31").Value = 10001
$ws.Range("B31").Value = 10030
$ws.Range("C31").Value = "eng"
$ws.Range("D31").Value = $true
$ws.Range("E31").Value = "superadmin"
$ws.Range("F31").Value = "now()"
$ws.Range("G31").Value = "now()"

# Row 32 - new Mac-Address entry
$ws.Range("A32").Value = 10001
$ws.Range("B32").Value = 10031
$ws.Range("C32").Value = "eng"
$ws.Range("D32").Value = $true
$ws.Range("E32").Value = "superadmin"
$ws.Range("F32").Value = "now()"
$ws.Range("G32").Value = "now()"

$ws.Range("E32").Select() | Out-Null
